$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 234, shifting existing rows 234-343 down to 235-344.
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with its data.
$ws.Cells.Item(234, 1).Value = 9
$ws.Cells.Item(234, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(234, 3).Value = "Metropolitana"
$ws.Cells.Item(234, 4).Value = 44636
$ws.Cells.Item(234, 4).NumberFormat = $ws.Cells.Item(235, 4).NumberFormat
$ws.Cells.Item(234, 5).Value = 13
$ws.Cells.Item(234, 6).Value = 100112032
$ws.Cells.Item(234, 7).Value = "Zapallo italiano"
$ws.Cells.Item(234, 8).Value = "Sin especificar"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 52
$ws.Cells.Item(234, 11).Value = 8000
$ws.Cells.Item(234, 12).Value = 9000
$ws.Cells.Item(234, 13).Value = 8500
$ws.Cells.Item(234, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(234, 15).Value = "Región Metropolitana"
$ws.Cells.Item(234, 16).Value = 142
$ws.Cells.Item(234, 17).Value = 60
$ws.Cells.Item(234, 18).Value = "Hortaliza"
